# Login.xlsx - "Correção de alguns promenores nos diagramas uml"
#
# The "Pós condição:" value (merged cell C5:D5) was corrected from
# "Ter acesso às suas informações" to "Autenticou-se no sistema", and the
# active selection on the sheet moved from D15 to D8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the post-condition text for the Login use case.
$ws.Range("C5").Value = "Autenticou-se no sistema"

# Restore the sheet's active cell/selection.
$ws.Range("D8").Select()
